$d = $word.ActiveDocument

# Locate the paragraph that contains "realizar a clusterização dos dados"
# (the last bullet under "O sistema deve permitir o upload de um arquivo .csv...").
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*realizar a clusteriza*dos dados*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    Write-Output "ERROR: target paragraph not found"
} else {
    $r = $target.Range

    # Rebuild the existing paragraph's OOXML exactly (preserving its run
    # boundaries/rsid) followed by the brand-new bullet requested in the
    # commit. InsertXML replaces the contents of the Range it is called on,
    # so calling it on the whole target paragraph's Range and re-supplying
    # that paragraph's own markup plus the new one appends the new bullet
    # right after it without disturbing anything else in the document.
    $origParaXml = '<w:p w14:paraId="124E8A54" w14:textId="0210BBD2" w:rsidR="00D65216" w:rsidRPr="007D2059" w:rsidRDefault="00D65216" w:rsidP="007D2059">' +
        '<w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
        '<w:r><w:t xml:space="preserve">O sistema deve permitir o </w:t></w:r>' +
        '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>upload</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> de um arquivo .csv e a partir dele </w:t></w:r>' +
        '<w:r w:rsidR="00D7635C"><w:t>realizar a clusterização dos dados</w:t></w:r>' +
        '</w:p>'

    $newParaXml = '<w:p>' +
        '<w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
        '<w:r><w:t>O sistema deve permitir a geração de gráficos para o usuário (mapa de calor de acordo com o número de propriedades da cidade...)</w:t></w:r>' +
        '</w:p>'

    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $origParaXml + $newParaXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($xmlFrag)
    Write-Output "Inserted new requirement bullet after the clusterization paragraph."
}
